$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New review comment strings (added to sharedStrings via cell writes) ---
$sAlarmStop  = "TSH: 24/01/2020 : The Alarm stop when press stop or 1 min passes "
$sNormalTone = "TSH: 24/01/2020 : The normal tone "
$sButtons = @"
TSH: 24/1/2020 : You can use the 3 buttons as follows:
For Watch and Alarms modes 
1- Button 1: Mode change "Watch , Alarm , stop watch"
2-Button 2 : Adjust/Select for hours minutes 
3- Button 3 : Up increment to set the hours or minutes
For Stop watch the change in :
- Button 2 : will be used to start and pause the stop watch 
- Button 3 : used to reset the stop watch 
You can off the alarm by Up increment if the alarm is on 
"@
$sAlarmOnly = "TSH: 24/01/2020 : no just for the alarm "

# --- Column I: widen to fit the new review/answer text ---
$ws.Columns.Item(9).ColumnWidth = 56.5

# --- Add the answer/review text to column I for each existing question row ---
$ws.Range("I2").Value = $sAlarmStop
$ws.Range("I3").Value = $sNormalTone

$ws.Range("I4").Value = $sButtons
$ws.Range("I4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 132

$ws.Range("I5").Value = $sAlarmOnly

# --- Update the view: select A6 and drop the previous scroll position ---
$ws.Range("A6").Select()
